# Update "想去人数" (want-to-go count) figures across the workbook sheets,
# matching the refreshed data snapshot generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 8188
$wsExhibition.Range("F5").Value = 32595
$wsExhibition.Range("F15").Value = 416
$wsExhibition.Range("F22").Value = 89
$wsExhibition.Range("F23").Value = 728
$wsExhibition.Range("F25").Value = 852
$wsExhibition.Range("F29").Value = 645

# Sheet "演出" (Performance)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F5").Value = 319
$wsPerformance.Range("F7").Value = 52

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 8188
$wsAll.Range("F7").Value = 32595
$wsAll.Range("F16").Value = 319
$wsAll.Range("F21").Value = 416
$wsAll.Range("F22").Value = 52
$wsAll.Range("F33").Value = 89
$wsAll.Range("F34").Value = 728
$wsAll.Range("F36").Value = 852
$wsAll.Range("F41").Value = 645
$wsAll.Range("F42").Value = 645
